$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 68, shifting existing rows 68-168 down to 69-169.
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with its data (matching the rest of the
# table's layout: A..R).
$ws.Range("A68").Value = 5
$ws.Range("B68").Value = "Macroferia Regional de Talca"
$ws.Range("C68").Value = "Maule"
$ws.Range("D68").Value = 44580
$ws.Range("E68").Value = 7
$ws.Range("F68").Value = 100112021
$ws.Range("G68").Value = "Ají"
$ws.Range("H68").Value = "Americana (o)"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 160
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("M68").Value = 15000
$ws.Range("N68").Value = "$/caja 14 kilos"
$ws.Range("O68").Value = "Región del Maule"
$ws.Range("P68").Value = 1071
$ws.Range("Q68").Value = 14
$ws.Range("R68").Value = "Hortaliza"
